$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new date header for 29-jun in column S
$ws.Range("S1").Value = "29-jun"

# Fill in the new column S values (one per data row)
$values = @(
    0,
    16.214590788615514,
    11.871641373282211,
    17.723916505413239,
    0,
    5.3872246234134087,
    5.4518904379050817,
    15.157706403255874,
    20.587123418261537,
    11.691203479300381,
    0,
    10.738703718410937,
    0,
    0,
    11.601964072329285,
    0,
    0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 19).Value = $values[$i]
}

# Update the active selection to match the post-edit state
$ws.Range("T7").Select()
